# Fill in the "Kemaro" custom case numbers on the Inp_Custom sheet, and
# point ModelInput!F1 at that sheet so the rest of the model (ModelOutput,
# charts, etc.) recalculates against these new assumptions.

$wb = $excel.ActiveWorkbook

$wsCustom = $wb.Worksheets.Item("Inp_Custom")
$wsModelInput = $wb.Worksheets.Item("ModelInput")
$wsModelOutput = $wb.Worksheets.Item("ModelOutput")

# --- Update the custom-case input values ---------------------------------
$wsCustom.Range("B8").Value = 2
$wsCustom.Range("B20").Value = 13000
$wsCustom.Range("B24").Value = 0
$wsCustom.Range("B30").Value = 26000
$wsCustom.Range("B32").Formula = "=1500/12"
$wsCustom.Range("B33").Value = 0
$wsCustom.Range("B41").Value = 0

# --- Switch the model to read from the Inp_Custom sheet -------------------
$wsModelInput.Range("F1").Value = "Inp_Custom"

# --- Update selections / scroll position to match the edited workbook -----
$wsCustom.Activate()
$wsCustom.Range("B9").Select()
$excel.ActiveWindow.ScrollRow = 3

$wsModelOutput.Activate()
$wsModelOutput.Range("B41").Select()
